# Apply the updates made to tests/A02_pixell_test_plan_investment_account.xlsx
# (new developer name, new Inputs/Setup/Expected-Result columns for the test
# plan rows, and the resulting cursor/selection move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Developer name (row 3) ---
$ws.Range("C3").Value = "Beerdavinder singh"

# --- Row 7: __init__ / Attributes are set to parameter values. ---
$ws.Range("E7").Value = "Inputs"
$ws.Range("F7").Value = "22222, 3333, 4444.44, date.today(), 4.00"
$ws.Range("G7").Value = "Setup"

# --- Row 8: __init__ / management fee has invalid type. ---
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = '22222, 3333, 4444.44, date.today(), "invalid"'
$ws.Range("G8").Value = "Raises ValueError"

# --- Row 9: __str__ / date created more than 10 years ago ---
$ws.Range("E9").Value = "(22222, 3333, 4444.44, date.today(), 4.00)"
$ws.Range("F9").Value = "None"
$ws.Range("G9").Value = "GETS SERVICE CHARGES"

# --- Row 10: __str__ / date created within last 10 years. ---
$ws.Range("E10").Value = "(22222, 3333, 4444.44, date.today(), 4.00)"
$ws.Range("F10").Value = "None"
$ws.Range("G10").Value = "GETS SERVICE CHARGES"

# --- Row 11: __str__ / date created exactly 10 years ago. ---
$ws.Range("E11").Value = "(22222, 3333, 4444.44, date.today(), 4.00)"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "GETS SERVICE CHARGES"

# --- Row 12: get_service_charges / displays waived management fee ... ---
$e12 = @'
 """
        self.investment = InvestmentAccount(22222, 3333, 4444.44, self.eleven_years_ago, 4.00)
        expected = f"Account number: 22222 Balance: $4444.44\n" \
                   f"Management Fee: Waived Account Type: Investment"
'@
$ws.Range("E12").Value = $e12
$ws.Range("F12").Value = "(22222, 3333, 4444.44, date.today(), 4.00)"
$ws.Range("G12").Value = "FORMAT STR"

# --- Row 13: get_service_charges / displays management fee ... ---
$e13 = @'
expected = f"Account number: 22222 Balance: $4444.44\n" \
                   f"Management Fee: $4.00 Account Type: Investment"
'@
$ws.Range("E13").Value = $e13
$ws.Range("F13").Value = "(22222, 3333, 4444.44, date.today(), 4.00)"
$ws.Range("G13").Value = "FORMAT STR"

# Final cursor/selection position, as last left by the author
$ws.Range("G13").Select()
